# Apply the 2026-02-05 18:44 JST append + re-sort update to the "ランサーズ" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove any existing hyperlinks; they will be re-created below at the correct rows
$ws.Hyperlinks.Delete()

# Row 2: Amazonから情報取得するツール作ってください。SP-AP...
$ws.Cells.Item(2,1).Value = "2026-02-05 18:44:13"
$ws.Cells.Item(2,2).Value = "Amazonから情報取得するツール作ってください。SP-API有 Python希望"
$ws.Cells.Item(2,3).Value = "システム開発"
$ws.Cells.Item(2,4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(2,5).Value = "期限情報なし"
$ws.Cells.Item(2,6).Value = "https://www.lancers.jp/work/detail/5486242"
$ws.Hyperlinks.Add($ws.Cells.Item(2,6), "https://www.lancers.jp/work/detail/5486242")
$ws.Cells.Item(2,6).Style = "Hyperlink"
$ws.Cells.Item(2,7).Value = 435
$ws.Cells.Item(2,8).Value = "🔥Python,API ◆ツール"

# Row 3: 産業機械向けAI異常検知・状態推定システムの開発・導入支援エ...
$ws.Cells.Item(3,1).Value = "2026-02-05 18:44:13"
$ws.Cells.Item(3,2).Value = "産業機械向けAI異常検知・状態推定システムの開発・導入支援エンジニア募集(AI/エッジ・組み込み)"
$ws.Cells.Item(3,3).Value = "システム開発"
$ws.Cells.Item(3,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(3,5).Value = "期限情報なし"
$ws.Cells.Item(3,6).Value = "https://www.lancers.jp/work/detail/5450864"
$ws.Hyperlinks.Add($ws.Cells.Item(3,6), "https://www.lancers.jp/work/detail/5450864")
$ws.Cells.Item(3,6).Style = "Hyperlink"
$ws.Cells.Item(3,7).Value = 383
$ws.Cells.Item(3,8).Value = "🔥AI,Ai ◆開発"

# Row 4: SaaSビジネスにおける「バーティカル(垂直型)」展開の横ス...
$ws.Cells.Item(4,1).Value = "2026-02-05 18:44:13"
$ws.Cells.Item(4,2).Value = "SaaSビジネスにおける「バーティカル(垂直型)」展開の横スライド可能なAIシステムの開発です"
$ws.Cells.Item(4,3).Value = "システム開発"
$ws.Cells.Item(4,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(4,5).Value = "期限情報なし"
$ws.Cells.Item(4,6).Value = "https://www.lancers.jp/work/detail/5485911"
$ws.Hyperlinks.Add($ws.Cells.Item(4,6), "https://www.lancers.jp/work/detail/5485911")
$ws.Cells.Item(4,6).Style = "Hyperlink"
$ws.Cells.Item(4,7).Value = 383
$ws.Cells.Item(4,8).Value = "🔥AI,Ai ◆開発"

# Row 5: 【フルタイム】最先端AI(LLM)開発エンジニア募集!新規プ...
$ws.Cells.Item(5,1).Value = "2026-02-05 18:44:13"
$ws.Cells.Item(5,2).Value = "【フルタイム】最先端AI(LLM)開発エンジニア募集!新規プロダクトの核となる開発パートナーを募集"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,6).Value = "https://www.lancers.jp/work/detail/5460294"
$ws.Hyperlinks.Add($ws.Cells.Item(5,6), "https://www.lancers.jp/work/detail/5460294")
$ws.Cells.Item(5,6).Style = "Hyperlink"
$ws.Cells.Item(5,7).Value = 375
$ws.Cells.Item(5,8).Value = "🔥AI,Ai ◆開発"

# Row 6: 【急募】n8nとTwitter APIを活用した自動化ワーク...
$ws.Cells.Item(6,1).Value = "2026-02-05 18:44:13"
$ws.Cells.Item(6,2).Value = "【急募】n8nとTwitter APIを活用した自動化ワークフロー作成依頼"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,6).Value = "https://www.lancers.jp/work/detail/5486501"
$ws.Hyperlinks.Add($ws.Cells.Item(6,6), "https://www.lancers.jp/work/detail/5486501")
$ws.Cells.Item(6,6).Style = "Hyperlink"
$ws.Cells.Item(6,7).Value = 255
$ws.Cells.Item(6,8).Value = "🔥API ◆自動化"

# Row 7: 【急募】クリックポスト自動発行ツール開発依頼...
$ws.Cells.Item(7,1).Value = "2026-02-05 18:44:13"
$ws.Cells.Item(7,2).Value = "【急募】クリックポスト自動発行ツール開発依頼"
$ws.Cells.Item(7,3).Value = "システム開発"
$ws.Cells.Item(7,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(7,5).Value = "期限情報なし"
$ws.Cells.Item(7,6).Value = "https://www.lancers.jp/work/detail/5485895"
$ws.Hyperlinks.Add($ws.Cells.Item(7,6), "https://www.lancers.jp/work/detail/5485895")
$ws.Cells.Item(7,6).Style = "Hyperlink"
$ws.Cells.Item(7,7).Value = 123
$ws.Cells.Item(7,8).Value = "◆ツール,開発"

# Row 8: 【Excelでのマクロ作成】リサーチツールの作成【スクレイピ...
$ws.Cells.Item(8,1).Value = "2026-02-05 18:44:13"
$ws.Cells.Item(8,2).Value = "【Excelでのマクロ作成】リサーチツールの作成【スクレイピング】"
$ws.Cells.Item(8,3).Value = "システム開発"
$ws.Cells.Item(8,4).Value = "1,000 ~ 5,000 円 / 固定"
$ws.Cells.Item(8,5).Value = "期限情報なし"
$ws.Cells.Item(8,6).Value = "https://www.lancers.jp/work/detail/5486225"
$ws.Hyperlinks.Add($ws.Cells.Item(8,6), "https://www.lancers.jp/work/detail/5486225")
$ws.Cells.Item(8,6).Style = "Hyperlink"
$ws.Cells.Item(8,7).Value = 100
$ws.Cells.Item(8,8).Value = "◆ツール,スクレイピング"

# Row 9: 初回 line予約システム、Googlrカレンダー連動一元管...
$ws.Cells.Item(9,1).Value = "2026-02-05 18:44:13"
$ws.Cells.Item(9,2).Value = "初回 line予約システム、Googlrカレンダー連動一元管理"
$ws.Cells.Item(9,3).Value = "システム開発"
$ws.Cells.Item(9,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(9,5).Value = "期限情報なし"
$ws.Cells.Item(9,6).Value = "https://www.lancers.jp/work/detail/5486342"
$ws.Hyperlinks.Add($ws.Cells.Item(9,6), "https://www.lancers.jp/work/detail/5486342")
$ws.Cells.Item(9,6).Style = "Hyperlink"
$ws.Cells.Item(9,7).Value = 53
$ws.Cells.Item(9,8).Value = "◇管理"

# Row 10: 【急募】iOS/AndroidアプリのSkyWay切替対応エ...
$ws.Cells.Item(10,1).Value = "2026-02-05 18:44:13"
$ws.Cells.Item(10,2).Value = "【急募】iOS/AndroidアプリのSkyWay切替対応エンジニア募集"
$ws.Cells.Item(10,3).Value = "システム開発"
$ws.Cells.Item(10,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(10,5).Value = "期限情報なし"
$ws.Cells.Item(10,6).Value = "https://www.lancers.jp/work/detail/5486110"
$ws.Hyperlinks.Add($ws.Cells.Item(10,6), "https://www.lancers.jp/work/detail/5486110")
$ws.Cells.Item(10,6).Style = "Hyperlink"
$ws.Cells.Item(10,7).Value = 38
$ws.Cells.Item(10,8).Value = "◇アプリ"

# Update column widths for columns B (title) and D (price) to fit the new content
$ws.Columns.Item(2).ColumnWidth = 50.166666666666664
$ws.Columns.Item(4).ColumnWidth = 31.166666666666664
